$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.166306018829346
$ws.Range("B1").Value = 2.437157869338989
$ws.Range("D1").Value = 2.36808443069458
$ws.Range("E1").Value = 1.234535813331604
